$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price correction on the existing "ISIS LS Gel 2,5L Lemon" line ---
$ws.Range("D28").Value2 = 490

# --- New product line ---
# Row 37 currently holds "Pril Isis Ultra Power 650ml" (item id 2970482, a
# text-typed / non-numeric-looking item code). That line is being pushed down
# to the bottom of the list (new row 39), and row 37 is being repurposed for a
# brand-new item: "LE CHAT power gel 4L" (id 2952095).
#
# Use Cut (move), not Copy/Paste, so the relocated row keeps both its original
# cell formatting (bordered table style) and its original value type (text)
# intact on row 39.
$ws.Range("A37:D37").Cut($ws.Range("A39:D39"))

# Row 37 is now free (still carrying the old bordered style) - fill it in with
# the new item. The id is a genuine number, unlike the row it replaced, and in
# the source data such numeric ids are stored borderless/unstyled, so clear
# the border that Cut left behind on A37.
$ws.Range("A37").Value2 = 2952095
$ws.Range("A37").Borders.LineStyle = -4142
$ws.Range("B37").Value2 = "LE CHAT power gel 4L"
$ws.Range("C37").Value2 = 3
$ws.Range("D37").Value2 = 1195

# --- Restore the scroll position / active selection recorded for this sheet ---
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("D29").Select()
